# Update the "F" column (collection/favorite count) values on the
# "展览" (Exhibitions) and "全部类型" (All types) worksheets to reflect the
# newly scraped counts.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibits) - first occurrence of the rows in the diff.
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value  = 117
$wsExhibit.Range("F9").Value  = 1302
$wsExhibit.Range("F10").Value = 279
$wsExhibit.Range("F12").Value = 10390
$wsExhibit.Range("F13").Value = 5
$wsExhibit.Range("F16").Value = 1021
$wsExhibit.Range("F18").Value = 11924
$wsExhibit.Range("F19").Value = 12309
$wsExhibit.Range("F20").Value = 29
$wsExhibit.Range("F21").Value = 111

# Sheet "全部类型" (all types) - second occurrence of the rows in the diff.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 117
$wsAll.Range("F10").Value = 1302
$wsAll.Range("F11").Value = 279
$wsAll.Range("F13").Value = 10391
$wsAll.Range("F14").Value = 5
$wsAll.Range("F17").Value = 1021
$wsAll.Range("F19").Value = 11924
$wsAll.Range("F20").Value = 12309
$wsAll.Range("F21").Value = 29
$wsAll.Range("F22").Value = 111

$wb.Save()
